# PnL_by_Day.xlsx — refresh of the "Query1" Power Query result table.
# The underlying query was re-run, pulling in two additional days of
# P&L data and shifting the remainder of the rows down (with a couple
# of rows re-sorted/re-valued to match the freshly refreshed source).
#
# This reproduces the resulting worksheet/table/defined-name state using
# the Excel COM object model (no direct OOXML surgery).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data (row, Pair, date-serial, P&L) after the refresh — rows 2..35.
$rows = @(
    @(2,  "BTCUSDT", 44612, -223.4),
    @(3,  "BTCUSDT", 44611, 41.31),
    @(4,  "BTCUSDT", 44610, 337.4),
    @(5,  "BTCUSDT", 44609, -166.98),
    @(6,  "BTCUSDT", 44608, -158.44999999999999),
    @(7,  "BTCUSDT", 44606, -103.18),
    @(8,  "BTCUSDT", 44605, -11.55),
    @(9,  "BTCUSDT", 44604, 260.14999999999998),
    @(10, "BTCUSDT", 44603, 252.52),
    @(11, "BTCUSDT", 44602, -244.53),
    @(12, "BTCUSDT", 44601, 5.25),
    @(13, "ETHUSDT", 44587, -290.8),
    @(14, "ETHUSDT", 44586, -831.64),
    @(15, "BTCUSDT", 44585, 14.26),
    @(16, "ETHUSDT", 44585, -45.64),
    @(17, "ETHUSDT", 44584, -387),
    @(18, "ETHUSDT", 44583, -874.83),
    @(19, "ETHUSDT", 44582, -1131.73),
    @(20, "ETHUSDT", 44581, 24.08),
    @(21, "BTCUSDT", 44580, 2.62),
    @(22, "ETHUSDT", 44580, -212.56),
    @(23, "BTCUSDT", 44579, -1832.65),
    @(24, "BTCUSDT", 44578, -104.96),
    @(25, "ETHUSDT", 44577, -30.73),
    @(26, "BTCUSDT", 44577, -103.03),
    @(27, "BTCUSDT", 44576, -85.29),
    @(28, "BTCUSDT", 44575, -18.149999999999999),
    @(29, "BTCUSDT", 44574, -43.46),
    @(30, "BTCUSDT", 44573, -1187.1300000000001),
    @(31, "BTCUSDT", 44572, -113.59),
    @(32, "BTCUSDT", 44571, -0.12),
    @(33, "ETHUSDT", 44296, -6.85),
    @(34, "BTCUSDT", 44280, -63.74),
    @(35, "BTCUSDT", 44279, 12.73)
)

# Reference cells that already carry the correct number formatting for
# the Pair (text) and date columns, used to stamp formatting onto the
# two brand-new rows (34 & 35) created by this refresh.
$fmtPairCell = $ws.Range("A2")
$fmtDateCell = $ws.Range("B2")

foreach ($r in $rows) {
    $rowNum = $r[0]
    $pair   = $r[1]
    $date   = $r[2]
    $pnl    = $r[3]

    if ($rowNum -gt 33) {
        # Brand-new rows added by the refresh: copy formatting down first.
        $fmtPairCell.Copy()
        $ws.Range("A$rowNum").PasteSpecial(-4122)
        $fmtDateCell.Copy()
        $ws.Range("B$rowNum").PasteSpecial(-4122)
    }

    $ws.Range("A$rowNum").Value = $pair
    $ws.Range("B$rowNum").Value = $date
    $ws.Range("C$rowNum").Value = $pnl
}

$excel.CutCopyMode = 0

# Grow the query table (ListObject) to cover the two extra rows.
$tbl = $ws.ListObjects.Item("Query1")
$tbl.Resize($ws.Range("A1:C35"))

# Keep the workbook-level "ExternalData_1" defined name (the query's
# external-data range) in sync with the new extent.
$wb.Names.Item("ExternalData_1").RefersTo = "=Sheet1!`$A`$1:`$C`$35"
